# Regen of save_data: column G ("K") is recomputed (std/mean + s_vals
# recalculated) instead of the old Strike# derived value. Update the
# 24 data rows (rows 2-25) with the newly computed K values.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("G2").Value = 5
$ws.Range("G3").Value = 5
$ws.Range("G4").Value = 5
$ws.Range("G5").Value = 6
$ws.Range("G6").Value = 3
$ws.Range("G7").Value = 4
$ws.Range("G8").Value = 7
$ws.Range("G9").Value = 3
$ws.Range("G10").Value = 7
$ws.Range("G11").Value = 2
$ws.Range("G12").Value = 3
$ws.Range("G13").Value = 5
$ws.Range("G14").Value = 3
$ws.Range("G15").Value = 5
$ws.Range("G16").Value = 3
$ws.Range("G17").Value = 3
$ws.Range("G18").Value = 4
$ws.Range("G19").Value = 5
$ws.Range("G20").Value = 2
$ws.Range("G21").Value = 6
$ws.Range("G22").Value = 3
$ws.Range("G23").Value = 4
$ws.Range("G24").Value = 4
$ws.Range("G25").Value = 1
